# Generate Report for Handback
# Updates handoff/handback timestamps for the eaf3e711-... file across the
# zh-cn and de-de language sheets, and refreshes the rolled-up
# "Latest HO Xliff Generate Date" on the Overview sheet accordingly.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# zh-cn: row 3 corresponds to eaf3e711-3d40-44c3-a4d4-9772a4a9983d.md
$wsZh.Range("H3").Value = "2016-08-21 06:53:56"
$wsZh.Range("K3").Value = "2016-08-21 06:54:24"

# de-de: row 3 corresponds to eaf3e711-3d40-44c3-a4d4-9772a4a9983d.md
$wsDe.Range("H3").Value = "2016-08-21 06:54:01"
$wsDe.Range("K3").Value = "2016-08-21 06:54:31"

# Overview: row 3 "Latest HO Xliff Generate Date" for eaf3e711-... reflects
# the newest handback datetime across languages (de-de handback: 06:54:31)
$wsOverview.Range("G3").Value = "2016-08-21 06:54:31"
